$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.523.02"
$ws.Range("E2").Value = "  +1.71%  "
$ws.Range("D3").Value = "1.825.54"
$ws.Range("E3").Value = "  +1.80%  "
$ws.Range("D4").Value = "'1.002"
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").Value = "'317.35"
$ws.Range("D7").Value = "'0.5398"
$ws.Range("E7").Value = "  +0.83%  "
$ws.Range("D8").Value = "'0.3993"
$ws.Range("E8").Value = "  +5.92%  "
$ws.Range("D9").Value = "'0.07731"
$ws.Range("E9").Value = "  +4.11%  "
$ws.Range("D10").Value = "'1.121"
$ws.Range("E10").Value = "  +2.63%  "
$ws.Range("D11").Value = "'42.03"
$ws.Range("E11").Value = "  +0.26%  "
$ws.Range("D12").Value = "'21.21"
$ws.Range("E12").Value = "  +3.12%  "
$ws.Range("D13").Value = "'6.348"
$ws.Range("E13").Value = "  +3.73%  "
$ws.Range("E14").Value = "  +0.17%  "
$ws.Range("D15").Value = "'7.609"
$ws.Range("E15").Value = "  +5.12%  "
$ws.Range("D16").Value = "1.826.67"
$ws.Range("E16").Value = "  +2.24%  "
$ws.Range("D17").Value = "'0.00001088"
$ws.Range("E17").Value = "  +2.85%  "
$ws.Range("D18").Value = "'89.90"
$ws.Range("E18").Value = "  +0.98%  "
$ws.Range("D19").Value = "'0.06583"
$ws.Range("E19").Value = "  +1.29%  "
$ws.Range("D20").Value = "'17.73"
$ws.Range("E20").Value = "  +2.78%  "
$ws.Range("E21").Value = "  +0.08%  "
$ws.Range("E22").Value = "  +3.03%  "
$ws.Range("D23").Value = "28.541.26"
$ws.Range("E23").Value = "  +1.69%  "
$ws.Range("D24").Value = "'11.20"
$ws.Range("E24").Value = "  +0.36%  "
$ws.Range("D25").Value = "'2.261"
$ws.Range("E25").Value = "  +8.18%  "
$ws.Range("D26").Value = "'158.17"
$ws.Range("E26").Value = "  +1.81%  "
$ws.Range("D27").Value = "'20.79"
$ws.Range("E27").Value = "  +2.42%  "
$ws.Range("D28").Value = "'2.452"
$ws.Range("E28").Value = "  +6.65%  "
$ws.Range("D29").Value = "2.039.28"
$ws.Range("E29").Value = "  +2.31%  "
$ws.Range("D30").Value = "'124.42"
$ws.Range("E30").Value = "  +2.72%  "
$ws.Range("E31").Value = "  +1.79%  "
$ws.Range("D32").Value = "'0.1125"
$ws.Range("E32").Value = "  +6.04%  "
$ws.Range("D33").Value = "'5.701"
$ws.Range("E33").Value = "  +2.60%  "
$ws.Range("B34").Value = "HuobiToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D34").Value = "'3.650"
$ws.Range("E34").Value = "  -0.36%  "
$ws.Range("B35").Value = "Hedera"
$ws.Range("C35").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D35").Value = "'0.07362"
$ws.Range("E35").Value = "  +13.62%  "
$ws.Range("D36").Value = "'0.2258"
$ws.Range("E36").Value = "  +0.53%  "
$ws.Range("D37").Value = "'0.02349"
$ws.Range("E37").Value = "  +2.65%  "
$ws.Range("D38").Value = "'8.945"
$ws.Range("E38").Value = "  +5.79%  "
$ws.Range("D39").Value = "'5.213"
$ws.Range("E39").Value = "  +3.85%  "
$ws.Range("D40").Value = "'11.40"
$ws.Range("E40").Value = "  +2.35%  "
$ws.Range("D41").Value = "'0.6293"
$ws.Range("E41").Value = "  +1.72%  "
$ws.Range("D42").Value = "'1.194"
$ws.Range("E42").Value = "  +1.29%  "
$ws.Range("E43").Value = "  +0.07%  "
$ws.Range("E44").Value = "  -3.33%  "
$ws.Range("D45").Value = "'13.54"
$ws.Range("E45").Value = "  +2.23%  "
$ws.Range("D46").Value = "'0.5902"
$ws.Range("E46").Value = "  +2.06%  "
$ws.Range("E47").Value = "  +1.13%  "
$ws.Range("D48").Value = "'125.31"
$ws.Range("E48").Value = "  +0.44%  "
$ws.Range("D49").Value = "'2.001"
$ws.Range("E49").Value = "  +3.85%  "
$ws.Range("E50").Value = "  +0.67%  "
$ws.Range("D51").Value = "'0.06923"
$ws.Range("E51").Value = "  +1.48%  "
